$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 45
$ws.Range("F6").Value = 1606
$ws.Range("F8").Value = 42
$ws.Range("F9").Value = 740
$ws.Range("F10").Value = 2696
$ws.Range("F11").Value = 2696
$ws.Range("F13").Value = 1788
$ws.Range("F14").Value = 613
$ws.Range("F15").Value = 288
$ws.Range("F16").Value = 701
$ws.Range("F17").Value = 5138
$ws.Range("F18").Value = 233
$ws.Range("F19").Value = 83
$ws.Range("F22").Value = 872
$ws.Range("F26").Value = 2441
$ws.Range("F28").Value = 373
$ws.Range("F31").Value = 488
$ws.Range("F32").Value = 1309
$ws.Range("F35").Value = 69
$ws.Range("F36").Value = 26
$ws.Range("F38").Value = 1460
$ws.Range("F39").Value = 22
$ws.Range("F40").Value = 1415

# Sheet: 演出 (Performance)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 116
$ws.Range("F11").Value = 150
$ws.Range("F15").Value = 27
$ws.Range("F16").Value = 143
$ws.Range("F17").Value = 335
$ws.Range("F18").Value = 262
$ws.Range("F19").Value = 519
$ws.Range("F25").Value = 5

# Sheet: 本地生活 (Local Life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 896
$ws.Range("F4").Value = 254
$ws.Range("F5").Value = 335
$ws.Range("F6").Value = 42
$ws.Range("F8").Value = 9

# Sheet: 全部类型 (All Types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 45
$ws.Range("F4").Value = 896
$ws.Range("F5").Value = 254
$ws.Range("F8").Value = 42
$ws.Range("F13").Value = 1606
$ws.Range("F16").Value = 42
$ws.Range("F17").Value = 2696
$ws.Range("F18").Value = 9
$ws.Range("F20").Value = 1788
$ws.Range("F21").Value = 150
$ws.Range("F22").Value = 613
$ws.Range("F23").Value = 288
$ws.Range("F24").Value = 701
$ws.Range("F25").Value = 5138
$ws.Range("F26").Value = 233
$ws.Range("F27").Value = 83
$ws.Range("F30").Value = 872
$ws.Range("F35").Value = 2441
$ws.Range("F37").Value = 373
$ws.Range("F39").Value = 488
$ws.Range("F40").Value = 1309
$ws.Range("F41").Value = 143
$ws.Range("F42").Value = 262
$ws.Range("F43").Value = 519
$ws.Range("F46").Value = 69
$ws.Range("F47").Value = 26
$ws.Range("F49").Value = 1415
